# Applies commit "Primary key of file table changed":
# Column B (file date / primary key) values are remapped to new dates,
# column widths for C and E are widened, and the active selection moves to B1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Staging column (far outside used range) used to round-trip cell formatting
# so that re-writing the cell text does not disturb its existing style/number format.
$ws.Range("B1").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("B1").Value = "'01/08/2023"
$ws.Range("Z2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("Z2").Clear()

$ws.Range("B2").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("B2").Value = "'01/08/2023"
$ws.Range("Z3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("Z3").Clear()

$ws.Range("B3").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("B3").Value = "'01/08/2023"
$ws.Range("Z4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("Z4").Clear()

$ws.Range("B4").Copy()
$ws.Range("Z5").PasteSpecial(-4122)
$ws.Range("B4").Value = "'01/08/2023"
$ws.Range("Z5").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("Z5").Clear()

$ws.Range("B5").Copy()
$ws.Range("Z6").PasteSpecial(-4122)
$ws.Range("B5").Value = "'01/08/2023"
$ws.Range("Z6").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("Z6").Clear()

$ws.Range("B6").Copy()
$ws.Range("Z7").PasteSpecial(-4122)
$ws.Range("B6").Value = "'01/08/2023"
$ws.Range("Z7").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("Z7").Clear()

$ws.Range("B7").Copy()
$ws.Range("Z8").PasteSpecial(-4122)
$ws.Range("B7").Value = "'01/08/2023"
$ws.Range("Z8").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("Z8").Clear()

$ws.Range("B8").Copy()
$ws.Range("Z9").PasteSpecial(-4122)
$ws.Range("B8").Value = "'01/08/2023"
$ws.Range("Z9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("Z9").Clear()

$ws.Range("B10").Copy()
$ws.Range("Z11").PasteSpecial(-4122)
$ws.Range("B10").Value = "'01/08/2023"
$ws.Range("Z11").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("Z11").Clear()

$ws.Range("B11").Copy()
$ws.Range("Z12").PasteSpecial(-4122)
$ws.Range("B11").Value = "'01/08/2023"
$ws.Range("Z12").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("Z12").Clear()

$ws.Range("B12").Copy()
$ws.Range("Z13").PasteSpecial(-4122)
$ws.Range("B12").Value = "'2023/08/01"
$ws.Range("Z13").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("Z13").Clear()

$ws.Range("B14").Copy()
$ws.Range("Z15").PasteSpecial(-4122)
$ws.Range("B14").Value = "'08/01/2023"
$ws.Range("Z15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("Z15").Clear()

$ws.Range("B16").Copy()
$ws.Range("Z17").PasteSpecial(-4122)
$ws.Range("B16").Value = "'01/08/2023"
$ws.Range("Z17").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("Z17").Clear()

$ws.Range("B17").Copy()
$ws.Range("Z18").PasteSpecial(-4122)
$ws.Range("B17").Value = "'01/08/2023"
$ws.Range("Z18").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("Z18").Clear()

# Widen columns C and E to fit the new, longer-looking date values
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 17

# Move the active selection to B1
[void]$ws.Range("B1").Select()
